# PM11 Tidsregistrering for Rasmus.xlsx - add a new time-registration row
# (row 7): a "Viderearbejde med mockup" task worked on 2020-02-25 from
# 08:40 to 10:30 under the "Graphic artist" role, logged as "2 timer".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 7 was a blank entry row (only number-formatted placeholder cells).
# Fill it in with the new task's data.
$ws.Range("A7").Value = "Viderearbejde med mockup"   # Opgavebeskrivelse
$ws.Range("B7").Value = "Graphic artist"             # Rolle
$ws.Range("C7").Value = 43886                        # Dato -> 2020-02-25
$ws.Range("D7").Value = 0.3611111111111111           # Starttid -> 08:40
$ws.Range("E7").Value = 0.4375                       # Sluttid -> 10:30
$ws.Range("F7").Value = "2 timer"                    # Aktuelt tidsforbrug

# Update the window/cursor position to reflect the new active cell and
# scroll position used when the workbook was saved.
$ws.Range("A8").Select()
$excel.ActiveWindow.ScrollRow = 2
